# Merge the split runs in the Title, Author and Abstract paragraphs into a
# single run each. The visible text is unchanged -- only how it is split
# across <w:r> runs collapses into one run per paragraph (matching the
# upstream commit). Find/Replace is scoped to each paragraph's own Range so
# that the identical text appearing later in the document (e.g. the "Zoë
# Gemmell" mention in the version-history section) is left untouched.

$d = $word.ActiveDocument

function Merge-ParagraphText($paragraph, $fullText) {
    $rng = $paragraph.Range
    $rng.Find.Execute(
        $fullText, $true, $false, $false, $false, $false,
        $true, 1, $false, $fullText, 2) | Out-Null
}

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Range.ParagraphStyle.NameLocal
    if ($styleName -eq "Title") {
        Merge-ParagraphText $p "Questions: Logarithms"
    } elseif ($styleName -eq "Author") {
        Merge-ParagraphText $p "Zoë Gemmell"
    } elseif ($styleName -eq "Abstract") {
        Merge-ParagraphText $p "A selection of questions for the study guide on logarithms."
    }
}
